$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Задание 1-2")

# Add the new shared string / label for the semi-perimeter row
$ws.Range("G11").Value = "Полупериметр треугольника"

# Semi-perimeter formula: (a+b+c)/2
$ws.Range("I11").Formula = "=(D4+D6+D8)/2"

# Fix the triangle-area formula to reuse the semi-perimeter cell (Heron's formula)
$ws.Range("I7").Formula = "=SQRT(I11*(I11-D4)*(I11-D6)*(I11-D8))"

# Update the active selection to match the post-edit cursor position
$ws.Range("N19").Select()
